# Clean the "Authors" column (E) data: the source text-join pipeline left
# accumulating runs of whitespace after each comma-separated author record,
# and one record ("Lei Xu") had a stale review-status flag. This mirrors the
# commit "cleaned more data, fixed a problem where some SLRs contained
# incorrect data."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Range("E" + $row)
    $value = $cell.Value2
    $cleaned = $value -replace ",( +)", ",  $1"
    $cell.Value = $cleaned
}

# Fix incorrect reviewer-status code for "Lei Xu" in the row 16 author list.
$fixCell = $ws.Range("E16")
$fixedValue = $fixCell.Value2 -replace "Lei%Xu%NULL%0,", "Lei%Xu%NULL%1,"
$fixCell.Value = $fixedValue
